$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 53, pushing the existing row 53 (and below) down to row 54.
$ws.Rows.Item(53).Insert()

# Populate the newly inserted row 53 with the new record's data.
$ws.Cells.Item(53, 1).Value = 1
$ws.Cells.Item(53, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(53, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(53, 4).Value = 44890
$ws.Cells.Item(53, 5).Value = 15
$ws.Cells.Item(53, 6).Value = "Fruta"
$ws.Cells.Item(53, 7).Value = 100108
$ws.Cells.Item(53, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(53, 9).Value = 100108001
$ws.Cells.Item(53, 10).Value = "Guayaba"
$ws.Cells.Item(53, 11).Value = "Sin especificar"
$ws.Cells.Item(53, 12).Value = "Segunda"
$ws.Cells.Item(53, 13).Value = 160
$ws.Cells.Item(53, 14).Value = 1000
$ws.Cells.Item(53, 15).Value = 1200
$ws.Cells.Item(53, 16).Value = 1100
$ws.Cells.Item(53, 17).Value = "$/kilo (en caja de 10 kilos )"
$ws.Cells.Item(53, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(53, 19).Value = 1100
$ws.Cells.Item(53, 20).Value = 1
